$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Day Wise Task Assigned / Status of Completion columns
# for the second team member's block (rows 6-9), which were left blank before.
$ws.Range("E6").Value = "User stories and step def creation"
$ws.Range("I6").Value = "Completed"

$ws.Range("E7").Value = "Payment page designing & Integration"
$ws.Range("I7").Value = "Completed"

$ws.Range("E8").Value = "User stories and step def creation"
$ws.Range("I8").Value = "Completed"

$ws.Range("E9").Value = "User stories and step def creation"
$ws.Range("I9").Value = "Completed"

# Update the active selection to match the author's final cursor position
$ws.Range("I9:M9").Select()
